$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.732.90"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.921.72"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.36"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.39"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.534"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.80"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "3.386.63"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.96"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "12.07"
$ws.Range("E15").Value = "  +63.40%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.55"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "2.971.70"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.985"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "50.696.81"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.02"
$ws.Range("E20").Value = "  -7.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.17"
$ws.Range("E21").Value = "  -5.27%  "
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.21"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.76"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  +8.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.79"
$ws.Range("E26").Value = "  -5.04%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.03"
$ws.Range("E28").Value = "  -5.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.29"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("E31").Value = "  -4.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.90"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.48"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.04"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.82"
$ws.Range("E35").Value = "  -4.90%  "
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.37"
$ws.Range("E40").Value = "  -4.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.42"
$ws.Range("E42").Value = "  -6.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.21"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.92"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.40"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "1.986.37"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.256"
$ws.Range("E49").Value = "  -6.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0317"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.23"
$ws.Range("E51").Value = "  +2.40%  "
